# Update NATMI TPM output values on the active worksheet.
# Columns: E=Ligand-expressing cells, F=Ligand detection rate,
#          G=Ligand average expression value, H=Ligand total expression value,
#          M=Receptor average expression value, N=Receptor total expression value,
#          O=Receptor derived specificity (avg), P=Receptor derived specificity (avg, dup),
#          Q=Edge average expression weight, R=Edge total expression weight,
#          S=Edge average expression derived specificity,
#          T=Edge total expression derived specificity

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.5039156666666667
$ws.Range("H2").Value = 1.511747
$ws.Range("M2").Value = 29.52617166666667
$ws.Range("N2").Value = 88.57851500000001
$ws.Range("O2").Value = 0.3218391660320701
$ws.Range("P2").Value = 0.3218391660320701
$ws.Range("Q2").Value = 14.87870047952278
$ws.Range("R2").Value = 133.908304315705
$ws.Range("S2").Value = 0.3218391660320701
$ws.Range("T2").Value = 0.3218391660320701

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.5039156666666667
$ws.Range("H3").Value = 1.511747
$ws.Range("O3").Value = 0.4328989896002822
$ws.Range("P3").Value = 0.4328989896002822
$ws.Range("Q3").Value = 20.013022291727
$ws.Range("R3").Value = 180.117200625543
$ws.Range("S3").Value = 0.4328989896002822
$ws.Range("T3").Value = 0.4328989896002822

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.5039156666666667
$ws.Range("H4").Value = 1.511747
$ws.Range("M4").Value = 22.50081433333333
$ws.Range("N4").Value = 67.502443
$ws.Range("O4").Value = 0.2452618443676477
$ws.Range("P4").Value = 0.2452618443676476
$ws.Range("Q4").Value = 11.33851285532456
$ws.Range("R4").Value = 102.046615697921
$ws.Range("S4").Value = 0.2452618443676477
$ws.Range("T4").Value = 0.2452618443676476
